$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.666.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.952.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.949.62'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.52%  '
$ws.Range("E11").Value = '  -4.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.461'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("E13").Value = '  -2.69%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.124'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.687.40'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.441.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.950.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.02%  '
$ws.Range("E20").Value = '  +13.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '445.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.698'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.37%  '
$ws.Range("E23").Value = '  -1.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("E25").Value = '  -2.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.52'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0000103'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.114'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.27'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.971'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.73'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '45.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.18'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.99'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.16%  '
$ws.Range("E41").Value = '  -2.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.121'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.56'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '386.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0352'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.681.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.85'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.99%  '
